$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$programsQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Special Topic",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details" 
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$projectsQuery = @'
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.project_org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;
'@

$grantsQuery = @'
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type  LIKE '%Cervical Cancer%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
'@

$publicationsQuery = @'
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
  WHEN pub.relative_citation_ratio = 1.0 THEN '1'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type  LIKE '%Cervical Cancer%'
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;
'@

$ws.Range("B2").Value = $programsQuery
$ws.Range("B3").Value = $projectsQuery
$ws.Range("B4").Value = $grantsQuery
$ws.Range("B5").Value = $publicationsQuery

$ws.Range("B2").Font.Size = 12
$ws.Range("C2").Font.Size = 12
$ws.Range("B3").Font.Size = 12
$ws.Range("B4").Font.Size = 12
$ws.Range("B5").Font.Size = 12
